$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '58.722.37'
$ws.Range('E2').Value2 = '  +2.23%  '
$ws.Range('D3').Value2 = '3.154.00'
$ws.Range('E3').Value2 = '  +2.45%  '
$ws.Range('E4').Value2 = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '534.73'
$ws.Range('E5').Value2 = '  +1.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '140.04'
$ws.Range('E6').Value2 = '  +2.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '1.00'
$ws.Range('E7').Value2 = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.515'
$ws.Range('E8').Value2 = '  +9.64%  '
$ws.Range('E9').Value2 = '  +1.43%  '
$ws.Range('E10').Value2 = '  +3.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '0.421'
$ws.Range('E11').Value2 = '  +4.09%  '
$ws.Range('E12').Value2 = '  +2.15%  '
$ws.Range('D13').Value2 = '3.699.73'
$ws.Range('E13').Value2 = '  +2.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '25.88'
$ws.Range('E14').Value2 = '  +2.60%  '
$ws.Range('E15').Value2 = '  +6.19%  '
$ws.Range('D16').Value2 = '58.789.69'
$ws.Range('E16').Value2 = '  +2.33%  '
$ws.Range('D17').Value2 = '3.157.49'
$ws.Range('E17').Value2 = '  +2.66%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '6.21'
$ws.Range('E18').Value2 = '  +5.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '13.02'
$ws.Range('E19').Value2 = '  +4.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '8.18'
$ws.Range('E20').Value2 = '  +4.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '372.33'
$ws.Range('E21').Value2 = '  +6.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '5.80'
$ws.Range('E22').Value2 = '  +2.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '0.999'
$ws.Range('E23').Value2 = '  +0.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '69.73'
$ws.Range('E24').Value2 = '  +2.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '0.513'
$ws.Range('E25').Value2 = '  +2.63%  '
$ws.Range('E26').Value2 = '  +0.94%  '
$ws.Range('E27').Value2 = '  +0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '8.03'
$ws.Range('E28').Value2 = '  +13.03%  '
$ws.Range('D29').Value2 = '0.0₃0874'
$ws.Range('E29').Value2 = '  +2.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '1.89'
$ws.Range('E30').Value2 = '  +2.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '6.16'
$ws.Range('E31').Value2 = '  +3.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '21.96'
$ws.Range('E32').Value2 = '  +4.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '5.19'
$ws.Range('E33').Value2 = '  +7.23%  '
$ws.Range('E34').Value2 = '  +3.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '159.81'
$ws.Range('E35').Value2 = '  +0.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '6.26'
$ws.Range('E36').Value2 = '  +4.43%  '
$ws.Range('E37').Value2 = '  +10.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '25.34'
$ws.Range('E38').Value2 = '  -0.33%  '
$ws.Range('B39').Value2 = 'Maker'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value2 = '2.654.19'
$ws.Range('E39').Value2 = '  +11.06%  '
$ws.Range('B40').Value2 = 'Stacks'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '1.68'
$ws.Range('E40').Value2 = '  +5.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '0.0684'
$ws.Range('E41').Value2 = '  +3.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '4.17'
$ws.Range('E42').Value2 = '  +4.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '38.82'
$ws.Range('E43').Value2 = '  +5.47%  '
$ws.Range('E44').Value2 = '  +2.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '0.0283'
$ws.Range('E45').Value2 = '  +8.83%  '
$ws.Range('E46').Value2 = '  +0.04%  '
$ws.Range('D47').Value2 = '3.197.62'
$ws.Range('E47').Value2 = '  +2.59%  '
$ws.Range('E48').Value2 = '  +12.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '0.986'
$ws.Range('E49').Value2 = '  +3.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '6.20'
$ws.Range('E50').Value2 = '  +3.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '20.32'
$ws.Range('E51').Value2 = '  +4.80%  '
